$wb = $excel.ActiveWorkbook

# --- classNumberOfLines: update "Number of Lines" (column B) ---
$wsClass = $wb.Worksheets.Item("classNumberOfLines")
$wsClass.Cells.Item(2, 2).NumberFormat = "@"
$wsClass.Cells.Item(2, 2).Value = "11"
$wsClass.Cells.Item(3, 2).NumberFormat = "@"
$wsClass.Cells.Item(3, 2).Value = "14"
$wsClass.Cells.Item(4, 2).NumberFormat = "@"
$wsClass.Cells.Item(4, 2).Value = "4"
$wsClass.Cells.Item(6, 2).NumberFormat = "@"
$wsClass.Cells.Item(6, 2).Value = "25"
$wsClass.Cells.Item(7, 2).NumberFormat = "@"
$wsClass.Cells.Item(7, 2).Value = "7"
$wsClass.Cells.Item(8, 2).NumberFormat = "@"
$wsClass.Cells.Item(8, 2).Value = "19"
$wsClass.Cells.Item(9, 2).NumberFormat = "@"
$wsClass.Cells.Item(9, 2).Value = "26"
$wsClass.Cells.Item(10, 2).NumberFormat = "@"
$wsClass.Cells.Item(10, 2).Value = "1"
$wsClass.Cells.Item(11, 2).NumberFormat = "@"
$wsClass.Cells.Item(11, 2).Value = "2"
$wsClass.Cells.Item(12, 2).NumberFormat = "@"
$wsClass.Cells.Item(12, 2).Value = "2"
$wsClass.Cells.Item(13, 2).NumberFormat = "@"
$wsClass.Cells.Item(13, 2).Value = "4"
$wsClass.Cells.Item(16, 2).NumberFormat = "@"
$wsClass.Cells.Item(16, 2).Value = "19"

# --- methodNumberOfLines: update "Number of Lines" (column C) ---
$wsMethod = $wb.Worksheets.Item("methodNumberOfLines")
$wsMethod.Cells.Item(2, 3).NumberFormat = "@"
$wsMethod.Cells.Item(2, 3).Value = "1"
$wsMethod.Cells.Item(3, 3).NumberFormat = "@"
$wsMethod.Cells.Item(3, 3).Value = "1"
$wsMethod.Cells.Item(5, 3).NumberFormat = "@"
$wsMethod.Cells.Item(5, 3).Value = "1"
$wsMethod.Cells.Item(6, 3).NumberFormat = "@"
$wsMethod.Cells.Item(6, 3).Value = "3"
$wsMethod.Cells.Item(8, 3).NumberFormat = "@"
$wsMethod.Cells.Item(8, 3).Value = "1"
$wsMethod.Cells.Item(9, 3).NumberFormat = "@"
$wsMethod.Cells.Item(9, 3).Value = "1"
$wsMethod.Cells.Item(10, 3).NumberFormat = "@"
$wsMethod.Cells.Item(10, 3).Value = "1"
$wsMethod.Cells.Item(11, 3).NumberFormat = "@"
$wsMethod.Cells.Item(11, 3).Value = "1"
$wsMethod.Cells.Item(12, 3).NumberFormat = "@"
$wsMethod.Cells.Item(12, 3).Value = "1"
$wsMethod.Cells.Item(14, 3).NumberFormat = "@"
$wsMethod.Cells.Item(14, 3).Value = "1"
$wsMethod.Cells.Item(15, 3).NumberFormat = "@"
$wsMethod.Cells.Item(15, 3).Value = "1"
$wsMethod.Cells.Item(16, 3).NumberFormat = "@"
$wsMethod.Cells.Item(16, 3).Value = "1"
$wsMethod.Cells.Item(18, 3).NumberFormat = "@"
$wsMethod.Cells.Item(18, 3).Value = "1"
$wsMethod.Cells.Item(19, 3).NumberFormat = "@"
$wsMethod.Cells.Item(19, 3).Value = "3"
$wsMethod.Cells.Item(22, 3).NumberFormat = "@"
$wsMethod.Cells.Item(22, 3).Value = "1"
$wsMethod.Cells.Item(23, 3).NumberFormat = "@"
$wsMethod.Cells.Item(23, 3).Value = "1"
$wsMethod.Cells.Item(26, 3).NumberFormat = "@"
$wsMethod.Cells.Item(26, 3).Value = "1"
$wsMethod.Cells.Item(27, 3).NumberFormat = "@"
$wsMethod.Cells.Item(27, 3).Value = "1"
$wsMethod.Cells.Item(28, 3).NumberFormat = "@"
$wsMethod.Cells.Item(28, 3).Value = "3"
$wsMethod.Cells.Item(30, 3).NumberFormat = "@"
$wsMethod.Cells.Item(30, 3).Value = "14"
$wsMethod.Cells.Item(32, 3).NumberFormat = "@"
$wsMethod.Cells.Item(32, 3).Value = "1"
$wsMethod.Cells.Item(33, 3).NumberFormat = "@"
$wsMethod.Cells.Item(33, 3).Value = "1"
$wsMethod.Cells.Item(34, 3).NumberFormat = "@"
$wsMethod.Cells.Item(34, 3).Value = "1"
$wsMethod.Cells.Item(35, 3).NumberFormat = "@"
$wsMethod.Cells.Item(35, 3).Value = "1"
$wsMethod.Cells.Item(36, 3).NumberFormat = "@"
$wsMethod.Cells.Item(36, 3).Value = "1"
$wsMethod.Cells.Item(37, 3).NumberFormat = "@"
$wsMethod.Cells.Item(37, 3).Value = "1"
$wsMethod.Cells.Item(42, 3).NumberFormat = "@"
$wsMethod.Cells.Item(42, 3).Value = "1"
$wsMethod.Cells.Item(43, 3).NumberFormat = "@"
$wsMethod.Cells.Item(43, 3).Value = "1"
$wsMethod.Cells.Item(44, 3).NumberFormat = "@"
$wsMethod.Cells.Item(44, 3).Value = "1"
$wsMethod.Cells.Item(45, 3).NumberFormat = "@"
$wsMethod.Cells.Item(45, 3).Value = "1"
$wsMethod.Cells.Item(50, 3).NumberFormat = "@"
$wsMethod.Cells.Item(50, 3).Value = "1"
$wsMethod.Cells.Item(51, 3).NumberFormat = "@"
$wsMethod.Cells.Item(51, 3).Value = "1"
$wsMethod.Cells.Item(52, 3).NumberFormat = "@"
$wsMethod.Cells.Item(52, 3).Value = "1"
$wsMethod.Cells.Item(54, 3).NumberFormat = "@"
$wsMethod.Cells.Item(54, 3).Value = "1"
$wsMethod.Cells.Item(55, 3).NumberFormat = "@"
$wsMethod.Cells.Item(55, 3).Value = "1"
$wsMethod.Cells.Item(57, 3).NumberFormat = "@"
$wsMethod.Cells.Item(57, 3).Value = "1"
$wsMethod.Cells.Item(58, 3).NumberFormat = "@"
$wsMethod.Cells.Item(58, 3).Value = "3"
$wsMethod.Cells.Item(60, 3).NumberFormat = "@"
$wsMethod.Cells.Item(60, 3).Value = "1"
